# Update the "meanEMG legmaxROM" values (first 4 data columns, B:E) on the
# single worksheet for both the header/count row (row 1) and the two data
# rows (CON / STR), then restore the selection to the updated range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - counts per group
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 - CON
$ws.Range("B2").Value = 18.503661238465703
$ws.Range("C2").Value = 4.7956491767371778
$ws.Range("D2").Value = 5.1091101460493578
$ws.Range("E2").Value = 1.1489913043664497

# Row 3 - STR
$ws.Range("B3").Value = 32.099311122121442
$ws.Range("C3").Value = 4.186635601313264
$ws.Range("D3").Value = -4.643452276585287
$ws.Range("E3").Value = 7.6170594777394296

# Match the saved selection from the source workbook.
$ws.Range("B1:E3").Select()
